$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation (dated 2022-09-06, serial 44810) was recorded for this
# weekly series. It belongs right after the existing row 33 and before the
# old row 34, so insert a fresh row at position 34 - this shifts the old
# rows 34..80 down to 35..81 (matching the new dimension A1:R81) - and then
# populate the new row with its data.
$ws.Rows("34:34").Insert()

$ws.Range("A34").Value = 10
$ws.Range("B34").Value = 'Vega Modelo de Temuco'
$ws.Range("C34").Value = 'La Araucanía'
$ws.Range("D34").Value = 44810
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 300000001
$ws.Range("G34").Value = 'Rabanito'
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 20
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 10000
$ws.Range("N34").Value = '$/docena de paquetes'
$ws.Range("O34").Value = 'Provincia de Cautín'
$ws.Range("P34").Value = 833
$ws.Range("Q34").Value = 12
$ws.Range("R34").Value = 'Hortaliza'
